# feat: add 2022-Q4 data
#
# 1) Insert a brand-new worksheet "2022-Q4" right before "2021-Q1" holding
#    the per-fund breakdown for the new quarter.
# 2) "总计" (sheet 1): insert a new row for "2022-Q4" right after the header,
#    pushing the existing "2021-Q1" / "2020-Q4" rows down one slot and
#    renumbering the A-column index (0,1,2,...).
#
# NOTE: worksheet handles in this host are position-based, not stable
# identities - once Worksheets.Add() reshuffles positions, any handle
# captured beforehand can silently start pointing at a different sheet.
# So the sheet-insert happens first, and every sheet used afterwards is
# re-fetched by name right before it's touched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" detail sheet before "2021-Q1"
# ---------------------------------------------------------------------
$q1sheet = $wb.Worksheets.Item("2021-Q1")
$ws = $wb.Worksheets.Add($q1sheet)
$ws.Name = "2022-Q4"

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# B (fund code) and D:G (decimal-looking figures) must stay text, so
# leading zeros ("010116") and trailing zeros ("0.3220") survive instead
# of being coerced to numbers.
$ws.Range("B2:B5").NumberFormat = "@"
$ws.Range("D2:G5").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "160322"
$ws.Range("C2").Value = "华夏港股通精选股票（LOF）A"
$ws.Range("D2").Value = "13.82"
$ws.Range("E2").Value = "92.59"
$ws.Range("F2").Value = "2.33"
$ws.Range("G2").Value = "0.3220"
$ws.Range("H2").Value = 9

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "010116"
$ws.Range("C3").Value = "民生加银新兴产业混合A"
$ws.Range("D3").Value = "6.99"
$ws.Range("E3").Value = "85.62"
$ws.Range("F3").Value = "3.74"
$ws.Range("G3").Value = "0.2614"
$ws.Range("H3").Value = 8

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "010117"
$ws.Range("C4").Value = "民生加银新兴产业混合C"
$ws.Range("D4").Value = "0.76"
$ws.Range("E4").Value = "85.62"
$ws.Range("F4").Value = "3.74"
$ws.Range("G4").Value = "0.0284"
$ws.Range("H4").Value = 8

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "012884"
$ws.Range("C5").Value = "华夏港股通精选股票（LOF）C"
$ws.Range("D5").Value = "0.69"
$ws.Range("E5").Value = "92.59"
$ws.Range("F5").Value = "2.33"
$ws.Range("G5").Value = "0.0161"
$ws.Range("H5").Value = 9

# Match the header / index-column formatting used by the sibling quarter
# sheets (bold, centered, thin border).
$q1fmt = $wb.Worksheets.Item("2021-Q1")
$q1fmt.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$q1fmt.Range("A2").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)

# Match the page margins used by the other quarter sheets.
$ps = $ws.PageSetup()
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push old row 3 ("2020-Q4") down to row 4, copying both value and format.
$total.Range("A4").Value = 2
$total.Range("B4").Value = $total.Range("B3").Value()
$total.Range("C4").Value = $total.Range("C3").Value()
$total.Range("D4").Value = $total.Range("D3").Value()
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

# Push old row 2 ("2021-Q1") down to row 3, copying both value and format.
$total.Range("A3").Value = 1
$total.Range("B3").Value = $total.Range("B2").Value()
$total.Range("C3").Value = $total.Range("C2").Value()
$total.Range("D3").Value = $total.Range("D2").Value()

# Write the new "2022-Q4" row into the now-free row 2.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.63

# ---------------------------------------------------------------------
# Restore the originally-active sheet ("总计", matching the unchanged
# bookViews/activeTab in the workbook part).
# ---------------------------------------------------------------------
$total2 = $wb.Worksheets.Item("总计")
$total2.Activate()

Write-Output "2022-Q4 sheet added"
